# Auto-generated edit script: updates "想去人数" (F) and "最低票价" (G) columns
# across sheets "展览" (1), "演出" (2) and "全部类型" (4)
# per gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("G3").Value = 160
$ws1.Range("F4").Value = 35
$ws1.Range("F5").Value = 135
$ws1.Range("F7").Value = 392
$ws1.Range("F8").Value = 4823
$ws1.Range("F9").Value = 4823
$ws1.Range("F10").Value = 21
$ws1.Range("F14").Value = 1104
$ws1.Range("F15").Value = 635
$ws1.Range("F16").Value = 4445
$ws1.Range("F17").Value = 176
$ws1.Range("F18").Value = 178
$ws1.Range("F19").Value = 77
$ws1.Range("F20").Value = 227
$ws1.Range("F21").Value = 3543
$ws1.Range("F22").Value = 5
$ws1.Range("F24").Value = 22
$ws1.Range("F25").Value = 3241
$ws1.Range("F26").Value = 142
$ws1.Range("F27").Value = 133
$ws1.Range("F30").Value = 203
$ws1.Range("F32").Value = 89
$ws1.Range("F33").Value = 68
$ws1.Range("F36").Value = 130
$ws1.Range("F37").Value = 5628
$ws1.Range("F38").Value = 880
$ws1.Range("F43").Value = 1139
$ws1.Range("F44").Value = 511
$ws1.Range("F46").Value = 2022
$ws1.Range("F47").Value = 303
$ws1.Range("F49").Value = 715
$ws1.Range("F50").Value = 864

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F8").Value = 114
$ws2.Range("F24").Value = 752

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("G3").Value = 160
$ws4.Range("F4").Value = 35
$ws4.Range("F6").Value = 135
$ws4.Range("F8").Value = 392
$ws4.Range("F9").Value = 4823
$ws4.Range("F10").Value = 4823
$ws4.Range("F11").Value = 21
$ws4.Range("F16").Value = 1104
$ws4.Range("F17").Value = 635
$ws4.Range("F18").Value = 4446
$ws4.Range("F19").Value = 176
$ws4.Range("F20").Value = 178
$ws4.Range("F21").Value = 77
$ws4.Range("F22").Value = 227
$ws4.Range("F23").Value = 3543
$ws4.Range("F24").Value = 3241
$ws4.Range("F25").Value = 142
$ws4.Range("F26").Value = 133
$ws4.Range("F28").Value = 203
$ws4.Range("F30").Value = 89
$ws4.Range("F31").Value = 68
$ws4.Range("F34").Value = 130
$ws4.Range("F36").Value = 5628
$ws4.Range("F38").Value = 880
$ws4.Range("F45").Value = 1139
$ws4.Range("F46").Value = 511
$ws4.Range("F47").Value = 2022
$ws4.Range("F48").Value = 303
$ws4.Range("F49").Value = 865

